$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add 8 new "work order" rows (29-36) below the existing data (rows 3-10),
#    mirroring the formatting of the last existing row (row 10) so the new
#    D/E (date) columns keep the same number format as the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A10:F10").Copy() | Out-Null
$ws.Range("A11:F18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Column A - orderId
$ws.Range("A11").Value = 29
$ws.Range("A12").Value = 30
$ws.Range("A13").Value = 31
$ws.Range("A14").Value = 32
$ws.Range("A15").Value = 33
$ws.Range("A16").Value = 34
$ws.Range("A17").Value = 35
$ws.Range("A18").Value = 36

# Column C - projectName (these introduce 8 brand new shared strings)
$ws.Range("C11").Value = "work order 29"
$ws.Range("C12").Value = "work order 30"
$ws.Range("C13").Value = "work order 31"
$ws.Range("C14").Value = "work order 32"
$ws.Range("C15").Value = "work order 33"
$ws.Range("C16").Value = "work order 34"
$ws.Range("C17").Value = "work order 35"
$ws.Range("C18").Value = "work order 36"

# Column E - shipDate (introduces 3 brand new shared strings, the rest reuse
# dates already present in the workbook)
$ws.Range("E12").Value = "18/03/2020 12:00PM"
$ws.Range("E14").Value = "10/06/2020 12:00PM"
$ws.Range("E11").Value = "30/03/2020 12:00PM"
$ws.Range("E13").Value = "05/05/2020 12:00PM"
$ws.Range("E15").Value = "08/08/2020 12:00PM"
$ws.Range("E16").Value = "05/05/2020 12:00PM"
$ws.Range("E17").Value = "10/07/2020 12:00PM"
$ws.Range("E18").Value = "08/08/2020 12:00PM"

# Column D - lastMaterialDate (reuses existing dates)
$ws.Range("D11").Value = "6/2/2020  12:00PM"
$ws.Range("D12").Value = "6/2/2020  12:00PM"
$ws.Range("D13").Value = "7/2/2020 12:00PM"
$ws.Range("D14").Value = "6/2/2020  12:00PM"
$ws.Range("D15").Value = "6/2/2020  12:00PM"
$ws.Range("D16").Value = "7/2/2020 12:00PM"
$ws.Range("D17").Value = "6/2/2020  12:00PM"
$ws.Range("D18").Value = "6/2/2020  12:00PM"

# Column B - partId
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 3
$ws.Range("B13").Value = 2
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 2
$ws.Range("B16").Value = 3
$ws.Range("B17").Value = 3
$ws.Range("B18").Value = 3

# Column F - quantity
$ws.Range("F11").Value = 4533
$ws.Range("F12").Value = 9583
$ws.Range("F13").Value = 2346
$ws.Range("F14").Value = 3582
$ws.Range("F15").Value = 2357
$ws.Range("F16").Value = 3572
$ws.Range("F17").Value = 2892
$ws.Range("F18").Value = 5652

# ---------------------------------------------------------------------------
# 2. Existing row 10 (orderId 28) had its quantity corrected.
# ---------------------------------------------------------------------------
$ws.Range("F10").Value = 2652

# ---------------------------------------------------------------------------
# 3. Update the active selection to reflect where the user ended up editing.
# ---------------------------------------------------------------------------
$null = $ws.Range("D7").Select()
